# Fluorescencia en el input y carpetas para P&D.
# Adds a "fluorescence" flag column (AD) to the Sheet1 input table, switches the
# example run label/comment from the old "TEST" row to a v8-without-fluorescence
# example (with new Id_min_run / Id_max_run values), renames the c*_star_660
# parameters to c_phy_star_660, and nudges a few column widths / the active
# selection to match the edited workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename A_c_star_660 / E_c_star_660 headers (N1/O1) to the *_phy_* variants.
$ws.Range("N1").Value = "A_c_phy_star_660"
$ws.Range("O1").Value = "E_c_phy_star_660"

# New column AD: fluorescence flag, with its header, example value and note.
$ws.Range("AD1").Value = "fluorescence"
$ws.Range("AD2").Value = 0
$ws.Range("AD3").Value = "0: sin, 1: con"

# Row 2 (example row): rename the label, update the comment, and add the
# min/max run id example values.
$ws.Range("A2").Value = "v8_no_fl"
$ws.Range("B2").Value = "Versión 8 – sin fluorescencia"
$ws.Range("C2").Value = 6000
$ws.Range("D2").Value = 8000

# Column width tweaks (B widened for the longer comment text, M/N/O resized
# to fit the renamed headers).
$ws.Range("B1").ColumnWidth = 23.92
$ws.Range("M1").ColumnWidth = 17.1
$ws.Range("N1").ColumnWidth = 15.7
$ws.Range("O1").ColumnWidth = 15.7

# Leave the cursor where the author left it when saving.
$ws.Range("I7").Select()
